$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = -0.06348608414014617
$ws.Range("D2").Value = 0.2150780536317457
$ws.Range("E2").Value = 0.004876747476953237
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -0.005803505760767809
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.19490196695336
$ws.Range("N2").Value = 0.0008140415976441889
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.08595621238231291
$ws.Range("V2").Value = 0.01660286692215323
$ws.Range("W2").Value = -0.005684954406440778
$ws.Range("Z2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0.03782319532730047
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.005604765158112685
$ws.Range("AF2").Value = -0.00216210829305749
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.02600752721072319
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.03476254392935594
$ws.Range("AO2").Value = 0.06772704628375234
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = 0
$ws.Range("AS2").Value = 0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1514092914727142
$ws.Range("AW2").Value = 0.0762025179529596
$ws.Range("AX2").Value = 0.008033511219956895
$ws.Range("AY2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = -0.01489799805321056
$ws.Range("BF2").Value = 0.06874910457888644
$ws.Range("BG2").Value = 0.02466433762295652
$ws.Range("BJ2").Value = 0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.0249503904145572
$ws.Range("BO2").Value = -0.04427272970730179
$ws.Range("BP2").Value = -0.07679347355527427
$ws.Range("BR2").Value = 0
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.06004265506089797
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.02387651955195982
$ws.Range("BY2").Value = -0.01591677721378175
$ws.Range("BZ2").Value = 0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = 0
$ws.Range("CE2").Value = 0.03130225729679532
$ws.Range("CF2").Value = 0
$ws.Range("CG2").Value = -0.03401870123622513
$ws.Range("CH2").Value = 0.0106228776687776
$ws.Range("CI2").Value = 0
$ws.Range("CJ2").Value = 0
$ws.Range("CL2").Value = 0
$ws.Range("CM2").Value = 0
$ws.Range("CN2").Value = -0.006928406815827574
$ws.Range("CO2").Value = 0
$ws.Range("CP2").Value = 0.02750862719667532
$ws.Range("CQ2").Value = 0.01950170115099826
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = 0
$ws.Range("CV2").Value = 0
$ws.Range("CW2").Value = 0.0439858022625931
$ws.Range("CX2").Value = 0
$ws.Range("CY2").Value = -0.0389323861523927
$ws.Range("CZ2").Value = 0.005248998491140817
$ws.Range("DE2").Value = 0
$ws.Range("DF2").Value = 0.02701288455505908
$ws.Range("DH2").Value = 0.02025060788531212
$ws.Range("DI2").Value = 0.0224068468429422
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = 0
$ws.Range("DL2").Value = 0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01737085642766599
$ws.Range("DQ2").Value = 0.0406539091754399
$ws.Range("DR2").Value = -0.004584347766881339
$ws.Range("DS2").Value = 0
$ws.Range("DU2").Value = 0
$ws.Range("DV2").Value = 0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.06290835086042507
$ws.Range("DY2").Value = 0
$ws.Range("DZ2").Value = -0.01570435231090992
$ws.Range("EA2").Value = -0.01193797399651102
$ws.Range("EB2").Value = 0
$ws.Range("ED2").Value = 0
$ws.Range("EF2").Value = 0
$ws.Range("EG2").Value = 0.03539537374972315
$ws.Range("EH2").Value = 0
$ws.Range("EI2").Value = 0.0815161155655162
$ws.Range("EJ2").Value = -0.03815334839236172
$ws.Range("EN2").Value = 0
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.0478223489034582
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03867461070842741
$ws.Range("ES2").Value = 0.01138094547169915
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.03762514914512365
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.03006370070916905
$ws.Range("FB2").Value = 0.009076250293154405
$ws.Range("FD2").Value = 0
$ws.Range("FF2").Value = 0
$ws.Range("FG2").Value = 0
$ws.Range("FH2").Value = -0.01811411271780501
$ws.Range("FJ2").Value = -0.007709830898137382
$ws.Range("FK2").Value = 0.02457764805577892
$ws.Range("FL2").Value = 0
$ws.Range("FP2").Value = 0
$ws.Range("FQ2").Value = -0.01100696790790678
$ws.Range("FR2").Value = 0
$ws.Range("FS2").Value = 0.002654848874162303
$ws.Range("FT2").Value = -0.007207252784700634
$ws.Range("FV2").Value = 0
$ws.Range("FW2").Value = 0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.02542218472942176
$ws.Range("GB2").Value = 0.01851740505980034
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = 0
